# "Finito l'esercizio aggiungendo gli attributi"
# Fill in the "Attributi (Constraints):" row (row 5) for every column that
# still needs a constraint: every attribute is NOT NULL, except the license
# plate ("Targa", column E) which is also UNIQUE since it identifies the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "NOT NULL"
$ws.Range("D5").Value = "NOT NULL"
$ws.Range("E5").Value = "NOT NULL, UNIQUE"
$ws.Range("F5").Value = "NOT NULL"
$ws.Range("G5").Value = "NOT NULL"
$ws.Range("H5").Value = "NOT NULL"
$ws.Range("I5").Value = "NOT NULL"
$ws.Range("J5").Value = "NOT NULL"
$ws.Range("K5").Value = "NOT NULL"
$ws.Range("L5").Value = "NOT NULL"
$ws.Range("M5").Value = "NOT NULL"
$ws.Range("N5").Value = "NOT NULL"
$ws.Range("O5").Value = "NOT NULL"
$ws.Range("P5").Value = "NOT NULL"

# The new, longer text in columns B (VIN) and E (Targa) no longer fits the
# previous column widths, so re-fit those two columns to their content.
$ws.Columns.Item(2).ColumnWidth = 24.333333333333332
$ws.Columns.Item(5).ColumnWidth = 16.666666666666668

# Leave the selection where the author ended up after typing the last value.
$ws.Range("H6").Select()
